$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.886.83"
$ws.Range("E2").Value = "  +1.69%  "

$ws.Range("D3").Value = "3.819.49"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'625.18"
$ws.Range("E5").Value = "  +4.12%  "

$ws.Range("D6").Value = "'165.11"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "3.819.08"
$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("E10").Value = "  +1.43%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("D12").Value = "'6.64"
$ws.Range("E12").Value = "  +4.04%  "

$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'35.88"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").Value = "4.457.84"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "3.872.99"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "68.878.81"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("D18").Value = "'18.16"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "'465.18"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").Value = "'9.68"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("D23").Value = "'0.706"
$ws.Range("E23").Value = "  +0.94%  "

$ws.Range("D24").Value = "'0.0000153"
$ws.Range("E24").Value = "  +4.96%  "

$ws.Range("D25").Value = "'83.85"
$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("D26").Value = "'12.02"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +2.29%  "

$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "3.969.80"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  -3.83%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.24"
$ws.Range("E32").Value = "  +2.30%  "

$ws.Range("D33").Value = "'7.34"
$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("D34").Value = "'29.15"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "'9.11"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("E37").Value = "  +1.95%  "

$ws.Range("E38").Value = "  +7.59%  "

$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +3.98%  "

$ws.Range("D40").Value = "'5.91"
$ws.Range("E40").Value = "  +2.81%  "

$ws.Range("D41").Value = "'0.978"
$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.301"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "'1.43"
$ws.Range("E45").Value = "  +2.87%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'154.69"
$ws.Range("E46").Value = "  +2.19%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'42.90"
$ws.Range("E47").Value = "  -5.14%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'46.80"
$ws.Range("E48").Value = "  -1.80%  "

$ws.Range("D49").Value = "'8.45"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("E50").Value = "  +2.33%  "

$ws.Range("D51").Value = "'380.98"
$ws.Range("E51").Value = "  -3.01%  "
